# edit.ps1
#
# Applies two changes described by the commit "Add intro to background
# section":
#
#   1. The "Clean Architecture..." heading run loses its explicit
#      <w:color w:val="3d3b49"/> run property.
#   2. Three new bibliography-style paragraphs are appended at the end
#      of the document (after filling in the previously-empty last
#      paragraph with the "Robert C Martin..." reference), introducing
#      the background section's additional references.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: drop the explicit color on the "Clean Architecture..." run.
# Rebuild just that paragraph (minus the bookmark, which is left alone
# so it isn't duplicated) with a minimal XML fragment so only this
# paragraph is touched.
# ---------------------------------------------------------------------
$headingRange = $d.Content
$headingRange.Find.Execute("Clean Architecture: A Craftsman") | Out-Null
$headingRange.Expand(4) | Out-Null   # wdParagraph -> grow to the full paragraph

$apos = [char]39
$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:between w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:fill="ffffff"/><w:spacing w:before="0" w:after="0" w:afterAutospacing="0" w:line="300" w:lineRule="auto"/><w:ind w:left="720" w:hanging="360"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Clean Architecture: A Craftsman' + $apos + 's Guide to Software Structure and Design, First Edition 2018 Pearson Education, Inc.</w:t></w:r></w:p>'
$headingRange.InsertXML($headingXml)

# ---------------------------------------------------------------------
# Step 2: turn the trailing empty paragraph into the "Robert C Martin"
# reference and append the new reference paragraphs after it.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> Robert C Martin. Clean code: a handbook of agile software craftsmanship. Pearson Education, 2009.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Igor Bugayenko. Elegant Objects. </w:t></w:r><w:r><w:rPr><w:color w:val="333333"/><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Createspace Independent Publishing Platform</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/><w:rPr><w:color w:val="333333"/><w:highlight w:val="white"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">J. R. Mckee, &#8220;Maintenance </w:t></w:r><w:r><w:rPr><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">as </w:t></w:r><w:r><w:rPr><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">a Function of Design&#8221;. Proceedings AFIPS, National Computer Conference, Las Vegas, pp 187-93. </w:t></w:r><w:r><w:rPr><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:color w:val="333333"/><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve"> </w:t><w:tab/><w:t xml:space="preserve"> </w:t><w:tab/><w:t xml:space="preserve"> </w:t><w:tab/><w:tab/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:rPr><w:color w:val="333333"/><w:highlight w:val="white"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="333333"/><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr></w:r></w:p>'

$lastPara.Range.InsertXML($newParagraphsXml)
